# Daily attendance processing - 2025-11-05 10:50:55
# For every row in the "Recorded By" column (G), swap the first two
# comma-separated entries (e.g. "System, user@x.com" -> "user@x.com, System").
# Entries with fewer than two comma-separated parts are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Value2

    if ($text -ne $null -and $text -ne "") {
        $parts = $text -split ", "
        if ($parts.Count -ge 2) {
            $first = $parts[0]
            $parts[0] = $parts[1]
            $parts[1] = $first
            $cell.Value2 = $parts -join ", "
        }
    }
}
